# Update PASADOR CERROJO price list: refresh the date and bump prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: date serial 45406 -> 45436 (2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# D29: 185.28 -> 364.992
$ws.Range("D29").Value = 364.992

# D30: 261.067 -> 514.29
$ws.Range("D30").Value = 514.29
